$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.025.51'
$ws.Range('E2').Value = '  -6.02%  '
$ws.Range('D3').Value = '3.262.97'
$ws.Range('E3').Value = '  -6.98%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '179.63'
$ws.Range('E5').Value = '  -11.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '517.98'
$ws.Range('E6').Value = '  -6.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.596'
$ws.Range('E7').Value = '  -0.70%  '
$ws.Range('D8').Value = '3.255.96'
$ws.Range('E8').Value = '  -6.99%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.617'
$ws.Range('E10').Value = '  -5.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.74'
$ws.Range('E11').Value = '  -5.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.131'
$ws.Range('E12').Value = '  -8.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  -6.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.11'
$ws.Range('E14').Value = '  -7.27%  '
$ws.Range('D15').Value = '3.764.20'
$ws.Range('E15').Value = '  -7.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.117'
$ws.Range('E16').Value = '  -5.81%  '
$ws.Range('D17').Value = '3.249.59'
$ws.Range('E17').Value = '  -7.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.55'
$ws.Range('E18').Value = '  -5.42%  '
$ws.Range('D19').Value = '62.948.97'
$ws.Range('E19').Value = '  -5.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.93'
$ws.Range('E20').Value = '  -7.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.947'
$ws.Range('E21').Value = '  -8.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '370.22'
$ws.Range('E22').Value = '  -4.83%  '
$ws.Range('E23').Value = '  -5.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.70'
$ws.Range('E24').Value = '  -7.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '79.65'
$ws.Range('E25').Value = '  -3.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.08'
$ws.Range('E26').Value = '  -0.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.82'
$ws.Range('E27').Value = '  +2.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.38'
$ws.Range('E28').Value = '  -5.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.62'
$ws.Range('E29').Value = '  -6.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.30'
$ws.Range('E30').Value = '  -6.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.41'
$ws.Range('E31').Value = '  -7.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.75'
$ws.Range('E32').Value = '  -7.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '636.06'
$ws.Range('E33').Value = '  -7.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.21'
$ws.Range('E34').Value = '  -4.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.106'
$ws.Range('E35').Value = '  -3.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.62'
$ws.Range('E36').Value = '  -7.29%  '
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.399'
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.30'
$ws.Range('E39').Value = '  -8.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.995'
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('D41').Value = '2.956.37'
$ws.Range('E41').Value = '  -5.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.124'
$ws.Range('E42').Value = '  -3.85%  '
$ws.Range('D43').Value = '0.0₃0654'
$ws.Range('E43').Value = '  -7.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.46'
$ws.Range('E44').Value = '  -2.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.68'
$ws.Range('E45').Value = '  -12.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0391'
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.82'
$ws.Range('E47').Value = '  +7.40%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.58'
$ws.Range('E48').Value = '  -5.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.125'
$ws.Range('E49').Value = '  -2.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.94'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.46'
$ws.Range('E51').Value = '  -12.69%  '
